$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws1.Range("H6").Value = 5138
$ws1.Range("I6").Value = 7240.143
$ws1.Range("J6").Value = 233
$ws1.Range("K6").Value = 21720.429
$ws1.Range("L6").Value = 699
$ws1.Range("M6").Value = -21608.429
$ws1.Range("N6").Value = -923
$ws1.Range("H8").Value = 381.125
$ws1.Range("I8").Value = 381.125
$ws1.Range("K8").Value = 1143.375
$ws1.Range("M8").Value = -1004.375
$ws1.Range("H40").Value = 2498.8572
$ws1.Range("I40").Value = 4767
$ws1.Range("J40").Value = 1880.2727
$ws1.Range("K40").Value = 4767
$ws1.Range("L40").Value = 1880.2727
$ws1.Range("M40").Value = -4592
$ws1.Range("N40").Value = -2230.2727
$ws1.Range("H52").Value = 800
$ws1.Range("I52").Value = 600
$ws1.Range("J52").Value = 900
$ws1.Range("K52").Value = 1800
$ws1.Range("L52").Value = 2700
$ws1.Range("M52").Value = -1640
$ws1.Range("N52").Value = -3020
$ws1.Range("H132").Value = 7941412.5
$ws1.Range("I132").Value = 8551990
$ws1.Range("K132").Value = 25655970
$ws1.Range("M132").Value = -25653440

$ws2 = $wb.Worksheets.Item("ARM")
$ws2.Range("H32").Value = 2978.7937
$ws2.Range("I32").Value = 2992.2642
$ws2.Range("J32").Value = 2907.4
$ws2.Range("K32").Value = 2992.2642
$ws2.Range("L32").Value = 2907.4
$ws2.Range("M32").Value = -2705.2642
$ws2.Range("N32").Value = -3481.4
$ws2.Range("H74").Value = 1844.4445
$ws2.Range("I74").Value = 757.25
$ws2.Range("J74").Value = 2714.2
$ws2.Range("K74").Value = 757.25
$ws2.Range("L74").Value = 2714.2
$ws2.Range("M74").Value = 116.75
$ws2.Range("N74").Value = -4462.2
$ws2.Range("H77").Value = 1844.4445
$ws2.Range("I77").Value = 757.25
$ws2.Range("J77").Value = 2714.2
$ws2.Range("K77").Value = 3786.25
$ws2.Range("L77").Value = 13571
$ws2.Range("M77").Value = 581.75
$ws2.Range("N77").Value = -22307
$ws2.Range("H110").Value = 1992.2
$ws2.Range("I110").Value = 1299
$ws2.Range("J110").Value = 2454.3333
$ws2.Range("K110").Value = 1299
$ws2.Range("L110").Value = 2454.3333
$ws2.Range("M110").Value = 746
$ws2.Range("N110").Value = -6544.3333
$ws2.Range("H132").Value = 2087.6572
$ws2.Range("I132").Value = 1798.3704
$ws2.Range("K132").Value = 5395.1112
$ws2.Range("M132").Value = -2865.1112

$ws3 = $wb.Worksheets.Item("BSM")
$ws3.Range("H86").Value = 3100.3125
$ws3.Range("I86").Value = 3057.0435
$ws3.Range("K86").Value = 3057.0435
$ws3.Range("M86").Value = -1934.0435
$ws3.Range("H89").Value = 3100.3125
$ws3.Range("I89").Value = 3057.0435
$ws3.Range("K89").Value = 15285.2175
$ws3.Range("M89").Value = -9669.217500000001
$ws3.Range("H105").Value = 166669840
$ws3.Range("I105").Value = 166669840
$ws3.Range("K105").Value = 166669840
$ws3.Range("M105").Value = -166668093
$ws3.Range("H107").Value = 2298.8333
$ws3.Range("I107").Value = 1620
$ws3.Range("J107").Value = 3656.5
$ws3.Range("K107").Value = 1620
$ws3.Range("L107").Value = 3656.5
$ws3.Range("M107").Value = 300
$ws3.Range("N107").Value = -7496.5
$ws3.Range("H113").Value = 3335733.2
$ws3.Range("I113").Value = 3335733.2
$ws3.Range("K113").Value = 3335733.2
$ws3.Range("M113").Value = -3333563.2
$ws3.Range("H130").Value = 30000
$ws3.Range("J130").Value = 30000
$ws3.Range("L130").Value = 30000
$ws3.Range("N130").Value = -40040

$ws4 = $wb.Worksheets.Item("CRP")
$ws4.Range("H31").Value = 1811.0435
$ws4.Range("I31").Value = 1360.8334
$ws4.Range("J31").Value = 2302.182
$ws4.Range("K31").Value = 1360.8334
$ws4.Range("L31").Value = 2302.182
$ws4.Range("M31").Value = -1065.8334
$ws4.Range("N31").Value = -2892.182
$ws4.Range("H34").Value = 1811.0435
$ws4.Range("I34").Value = 1360.8334
$ws4.Range("J34").Value = 2302.182
$ws4.Range("K34").Value = 1360.8334
$ws4.Range("L34").Value = 2302.182
$ws4.Range("M34").Value = -1158.8334
$ws4.Range("N34").Value = -2706.182
$ws4.Range("H58").Value = 1319.2778
$ws4.Range("J58").Value = 1751.5
$ws4.Range("L58").Value = 1751.5
$ws4.Range("N58").Value = -2157.5
$ws4.Range("H132").Value = 6519.875
$ws4.Range("I132").Value = 8959.786
$ws4.Range("J132").Value = 3104
$ws4.Range("K132").Value = 26879.358
$ws4.Range("L132").Value = 9312
$ws4.Range("M132").Value = -24349.358
$ws4.Range("N132").Value = -14372
$ws4.Range("H136").Value = 1319.2778
$ws4.Range("J136").Value = 1751.5
$ws4.Range("L136").Value = 5254.5
$ws4.Range("N136").Value = -10354.5

$ws5 = $wb.Worksheets.Item("CUL")
$ws5.Range("H5").Value = 1329.8276
$ws5.Range("I5").Value = 1427.9584
$ws5.Range("J5").Value = 858.8
$ws5.Range("K5").Value = 4283.8752
$ws5.Range("L5").Value = 2576.4
$ws5.Range("M5").Value = -4171.8752
$ws5.Range("N5").Value = -2800.4
$ws5.Range("H107").Value = 14814.143
$ws5.Range("I107").Value = 0
$ws5.Range("J107").Value = 14814.143
$ws5.Range("K107").Value = 0
$ws5.Range("L107").Value = 44442.429
$ws5.Range("M107").ClearContents()
$ws5.Range("N107").Value = -48282.429
$ws5.Range("H122").Value = 808.06665
$ws5.Range("I122").Value = 605.1667
$ws5.Range("J122").Value = 943.3333
$ws5.Range("K122").Value = 5446.5003
$ws5.Range("L122").Value = 8489.9997
$ws5.Range("M122").Value = -2996.5003
$ws5.Range("N122").Value = -13389.9997
$ws5.Range("H131").Value = 14286952
$ws5.Range("I131").Value = 250000260
$ws5.Range("J131").Value = 1297.3182
$ws5.Range("K131").Value = 750000780
$ws5.Range("L131").Value = 3891.9546
$ws5.Range("M131").Value = -749995740
$ws5.Range("N131").Value = -13971.9546
$ws5.Range("H135").Value = 1329.8276
$ws5.Range("I135").Value = 1427.9584
$ws5.Range("J135").Value = 858.8
$ws5.Range("K135").Value = 12851.6256
$ws5.Range("L135").Value = 7729.2
$ws5.Range("M135").Value = -10316.6256
$ws5.Range("N135").Value = -12799.2
$ws5.Range("H137").Value = 2144.6924
$ws5.Range("I137").Value = 1168
$ws5.Range("J137").Value = 7516.5
$ws5.Range("K137").Value = 3504
$ws5.Range("L137").Value = 22549.5
$ws5.Range("M137").Value = 1596
$ws5.Range("N137").Value = -32749.5

$ws6 = $wb.Worksheets.Item("GSM")
$ws6.Range("H80").Value = 6790
$ws6.Range("I80").Value = 6790
$ws6.Range("J80").Value = 0
$ws6.Range("K80").Value = 6790
$ws6.Range("L80").Value = 0
$ws6.Range("M80").Value = -5792
$ws6.Range("N80").ClearContents()
$ws6.Range("H83").Value = 6790
$ws6.Range("I83").Value = 6790
$ws6.Range("J83").Value = 0
$ws6.Range("K83").Value = 33950
$ws6.Range("L83").Value = 0
$ws6.Range("M83").Value = -28958
$ws6.Range("N83").ClearContents()
$ws6.Range("H132").Value = 2076.2415
$ws6.Range("I132").Value = 1857.762
$ws6.Range("K132").Value = 5573.286
$ws6.Range("M132").Value = -3043.286

$ws7 = $wb.Worksheets.Item("LTW")
$ws7.Range("H22").Value = 901.4737
$ws7.Range("J22").Value = 1331.5555
$ws7.Range("L22").Value = 1331.5555
$ws7.Range("N22").Value = -1921.5555
$ws7.Range("H27").Value = 901.4737
$ws7.Range("J27").Value = 1331.5555
$ws7.Range("L27").Value = 1331.5555
$ws7.Range("N27").Value = -1545.5555
$ws7.Range("H40").Value = 2931.182
$ws7.Range("I40").Value = 2700.8
$ws7.Range("J40").Value = 3123.1667
$ws7.Range("K40").Value = 2700.8
$ws7.Range("L40").Value = 3123.1667
$ws7.Range("M40").Value = -2564.8
$ws7.Range("N40").Value = -3395.1667
$ws7.Range("H55").Value = 238.85185
$ws7.Range("I55").Value = 158.25
$ws7.Range("J55").Value = 303.33334
$ws7.Range("K55").Value = 158.25
$ws7.Range("L55").Value = 303.33334
$ws7.Range("M55").Value = 14.75
$ws7.Range("N55").Value = -649.33334
$ws7.Range("H68").Value = 2011.6666
$ws7.Range("I68").Value = 1801.1666
$ws7.Range("K68").Value = 1801.1666
$ws7.Range("M68").Value = -1052.1666
$ws7.Range("H71").Value = 2011.6666
$ws7.Range("I71").Value = 1801.1666
$ws7.Range("K71").Value = 9005.833000000001
$ws7.Range("M71").Value = -5261.833000000001
$ws7.Range("H82").Value = 1425.4546
$ws7.Range("I82").Value = 1372.5
$ws7.Range("J82").Value = 1566.6666
$ws7.Range("K82").Value = 1372.5
$ws7.Range("L82").Value = 1566.6666
$ws7.Range("M82").Value = -1011.5
$ws7.Range("N82").Value = -2288.6666
$ws7.Range("H85").Value = 1425.4546
$ws7.Range("I85").Value = 1372.5
$ws7.Range("J85").Value = 1566.6666
$ws7.Range("K85").Value = 1372.5
$ws7.Range("L85").Value = 1566.6666
$ws7.Range("M85").Value = -124.5
$ws7.Range("N85").Value = -4062.6666
$ws7.Range("H127").Value = 0
$ws7.Range("J127").Value = 0
$ws7.Range("L127").Value = 0
$ws7.Range("N127").ClearContents()
$ws7.Range("H132").Value = 19829.2
$ws7.Range("I132").Value = 1241.3715
$ws7.Range("J132").Value = 52357.9
$ws7.Range("K132").Value = 3724.1145
$ws7.Range("L132").Value = 157073.7
$ws7.Range("M132").Value = -1194.1145
$ws7.Range("N132").Value = -162133.7

$ws8 = $wb.Worksheets.Item("WVR")
$ws8.Range("H81").Value = 505
$ws8.Range("J81").Value = 0
$ws8.Range("L81").Value = 0
$ws8.Range("N81").ClearContents()
$ws8.Range("H84").Value = 505
$ws8.Range("J84").Value = 0
$ws8.Range("L84").Value = 0
$ws8.Range("N84").ClearContents()
$ws8.Range("H132").Value = 3330.8333
$ws8.Range("I132").Value = 4799.8
$ws8.Range("J132").Value = 2765.8462
$ws8.Range("K132").Value = 14399.4
$ws8.Range("L132").Value = 8297.5386
$ws8.Range("M132").Value = -13357.5386

Write-Host "Applied all changes"